# Applies the commit "renamed and consolidated some parts names" to the
# Nestbox BOM workbook:
#   - drops the now-resolved note "consolidate with other inductor
#     possible?" from L201's NOTES cell (M34)
#   - marks the remaining NOTES-column remarks as addressed by recoloring
#     their text green (theme accent3) instead of leaving the yellow-filled
#     default look; the C602-position remark (M17) additionally shrinks to
#     11pt
#   - unhides/resizes the previously-hidden MFG/MPN/MFR columns (F:L) so the
#     part-number consolidation is visible
#   - updates the saved view (zoom + scroll position) to where the edits
#     were made

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- NOTES column: remove the resolved inductor-consolidation remark ---
$ws.Range("M34").ClearContents()

# --- NOTES column: recolor the still-open remarks green (theme accent3) ---
$greenCells = @("M8","M14","M16","M19","M35","M45","M46","M51","M54","M59","M63")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Font.ThemeColor = 7
}

# M17 (C602 label note) gets the same green but a smaller 11pt font
$ws.Range("M17").Font.ThemeColor = 7
$ws.Range("M17").Font.Size = 11

# --- Unhide & resize the part-number columns (F:L) that were hidden ---
$ws.Range("F1:G1").EntireColumn.Hidden = $false
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 25.5

$ws.Range("H1").EntireColumn.Hidden = $false
$ws.Range("H1").EntireColumn.ColumnWidth = 7

$ws.Range("I1:L1").EntireColumn.Hidden = $false
$ws.Range("I1:L1").EntireColumn.ColumnWidth = 7.83

# --- Update the saved view: zoom + scroll/selection position ---
$ws.Range("D57").Select()
$aw = $excel.ActiveWindow
$aw.Zoom = 125
$aw.ScrollRow = 36
$aw.ScrollColumn = 2
